$wb = $excel.ActiveWorkbook

$cur    = $wb.Worksheets.Item("Current")
$rg     = $wb.Worksheets.Item("RG table")
$season = $wb.Worksheets.Item("Season Log")

# -----------------------------------------------------------------
# Sheet: Season Log - insert the new Draftshot result for the
# Baltimore Orioles righties stack above the (old) row 134.
# This is done first so the new shared string for that stack name
# is created before the other new strings below (matches the order
# new strings end up appended to the shared string table).
# -----------------------------------------------------------------
$season.Rows(134).Insert()
$season.Range("A134").Value = 43610
$season.Range("B134").Value = "Draftshot"
$season.Range("C134").Value = "Baltimore Orioles righties (FD, DK)"
$season.Range("D134").Value = 5.39
$season.Range("E134").Value = "Success"

# -----------------------------------------------------------------
# Sheet: Current
# -----------------------------------------------------------------

# First stack (A:D) - Baltimore Orioles hitters
$cur.Range("A1").Value = "Baltimore Orioles hitters (FD, DK)"

$cur.Range("A3").Value = "Alberto"
$cur.Range("B3").Value = 2700
$cur.Range("C3").Value = 3

$cur.Range("A4").Value = "Nunez"
$cur.Range("B4").Value = 3500
$cur.Range("C4").Value = 15.5

$cur.Range("A6").Value = "Broxton"
$cur.Range("B6").Value = 2300
$cur.Range("C6").Value = 18.7

$cur.Range("A5").Value = "Stewart"
$cur.Range("B5").Value = 2300
$cur.Range("C5").Value = 0

# Second stack (F:I) - Arizona Diamondbacks hitters
$cur.Range("F1").Value = "Arizona Diamondbacks hitters (FD, DK)"

$cur.Range("F3").Value = "Locastro"
$cur.Range("G3").Value = 2500
$cur.Range("H3").Value = 0

$cur.Range("F4").Value = "Jones"
$cur.Range("G4").Value = 4000
$cur.Range("H4").Value = 18.7

$cur.Range("F5").Value = "Cron"
$cur.Range("G5").Value = 2800
$cur.Range("H5").Value = 15.2

$cur.Range("F6").Value = "Vargas"
$cur.Range("G6").Value = 3100
$cur.Range("H6").Value = 15

# Third stack (K:N) is no longer tracked on this sheet - clear it out
$cur.Range("K1").ClearContents()
$cur.Range("K3:K6").ClearContents()
$cur.Range("L3:L6").ClearContents()
$cur.Range("M3:M6").ClearContents()

# Update view selection on Current sheet
$cur.Activate()
$cur.Range("F1:I1").Select()

# -----------------------------------------------------------------
# Sheet: RG table
# -----------------------------------------------------------------
$rg.Activate()
$rg.Range("A1:A2").Select()

# -----------------------------------------------------------------
# Sheet: Season Log (remaining additions)
# -----------------------------------------------------------------

# Append the two new RG results (the stacks updated above on Current sheet)
$season.Range("A140").Value = 43614
$season.Range("B140").Value = "RG"
$season.Range("C140").Value = "Baltimore Orioles hitters (FD, DK)"
$season.Range("D140").Value = 3.44
$season.Range("E140").Value = "Failure"

$season.Range("A141").Value = 43614
$season.Range("B141").Value = "RG"
$season.Range("C141").Value = "Arizona Diamondbacks hitters (FD, DK)"
$season.Range("D141").Value = 3.94
$season.Range("E141").Value = "Success"

# Start of a new (not yet complete) Draftshot entry
$season.Range("A142").Value = 43615
$season.Range("B142").Value = "Draftshot"

# Update view selection / scroll on Season Log sheet
$season.Activate()
$excel.ActiveWindow.ScrollRow = 125
$season.Range("C142").Select()

$wb.Save()
